$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 548:549, shifting the existing rows (old 548..629) down to 550..631.
$ws.Rows("548:549").Insert()

# Fill in the new row 548 with its data.
$ws.Cells.Item(548, 1).Value = 10
$ws.Cells.Item(548, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(548, 3).Value = 'La Araucanía'
$ws.Cells.Item(548, 4).Value = 45015
$ws.Cells.Item(548, 5).Value = 9
$ws.Cells.Item(548, 6).Value = 100112023
$ws.Cells.Item(548, 7).Value = 'Brócoli'
$ws.Cells.Item(548, 8).Value = 'Sin especificar'
$ws.Cells.Item(548, 9).Value = 'Primera'
$ws.Cells.Item(548, 10).Value = 500
$ws.Cells.Item(548, 11).Value = 1300
$ws.Cells.Item(548, 12).Value = 1400
$ws.Cells.Item(548, 13).Value = 1340
$ws.Cells.Item(548, 14).Value = '$/unidad'
$ws.Cells.Item(548, 15).Value = 'Provincia de Cautín'
$ws.Cells.Item(548, 16).Value = 1340
$ws.Cells.Item(548, 17).Value = 1
$ws.Cells.Item(548, 18).Value = 'Hortaliza'

# Fill in the new row 549 with its data.
$ws.Cells.Item(549, 1).Value = 10
$ws.Cells.Item(549, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(549, 3).Value = 'La Araucanía'
$ws.Cells.Item(549, 4).Value = 45015
$ws.Cells.Item(549, 5).Value = 9
$ws.Cells.Item(549, 6).Value = 100112023
$ws.Cells.Item(549, 7).Value = 'Brócoli'
$ws.Cells.Item(549, 8).Value = 'Sin especificar'
$ws.Cells.Item(549, 9).Value = 'Primera'
$ws.Cells.Item(549, 10).Value = 2800
$ws.Cells.Item(549, 11).Value = 1300
$ws.Cells.Item(549, 12).Value = 1300
$ws.Cells.Item(549, 13).Value = 1300
$ws.Cells.Item(549, 14).Value = '$/unidad'
$ws.Cells.Item(549, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(549, 16).Value = 1300
$ws.Cells.Item(549, 17).Value = 1
$ws.Cells.Item(549, 18).Value = 'Hortaliza'
